# Update NATMI TPM-derived metrics on the active worksheet to reflect
# results recomputed with the new TPM values (see commit: "update scripts
# wuth new tpm"). Only numeric result columns (G..T, excluding the
# unchanged count columns) are affected; identifiers in A:F/K/L are
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5277743333333333
$ws.Range("H2").Value = 1.583323
$ws.Range("I2").Value = 0.01400965908295571
$ws.Range("J2").Value = 0.01400965908295571
$ws.Range("M2").Value = 97.57717366666668
$ws.Range("N2").Value = 292.731521
$ws.Range("O2").Value = 0.3532166605548384
$ws.Range("P2").Value = 0.3532166605548384
$ws.Range("Q2").Value = 51.4987277804759
$ws.Range("R2").Value = 463.4885500242831
$ws.Range("S2").Value = 0.004948444996793375
$ws.Range("T2").Value = 0.004948444996793375

# Row 3
$ws.Range("G3").Value = 0.5277743333333333
$ws.Range("H3").Value = 1.583323
$ws.Range("I3").Value = 0.01400965908295571
$ws.Range("J3").Value = 0.01400965908295571
$ws.Range("O3").Value = 0.5533024543641269
$ws.Range("P3").Value = 0.5533024543641269
$ws.Range("Q3").Value = 80.67108848378767
$ws.Range("R3").Value = 726.039796354089
$ws.Range("S3").Value = 0.007751578755404078
$ws.Range("T3").Value = 0.007751578755404076

# Row 4
$ws.Range("G4").Value = 0.5277743333333333
$ws.Range("H4").Value = 1.583323
$ws.Range("I4").Value = 0.01400965908295571
$ws.Range("J4").Value = 0.01400965908295571
$ws.Range("O4").Value = 0.09348088508103472
$ws.Range("P4").Value = 0.09348088508103473
$ws.Range("Q4").Value = 13.629443882698
$ws.Range("R4").Value = 122.664994944282
$ws.Range("S4").Value = 0.001309635330758257
$ws.Range("T4").Value = 0.001309635330758257

# Row 5
$ws.Range("I5").Value = 0.9003255417707673
$ws.Range("J5").Value = 0.9003255417707672
$ws.Range("M5").Value = 97.57717366666668
$ws.Range("N5").Value = 292.731521
$ws.Range("O5").Value = 0.3532166605548384
$ws.Range("P5").Value = 0.3532166605548384
$ws.Range("Q5").Value = 3309.546628859164
$ws.Range("R5").Value = 29785.91965973247
$ws.Range("S5").Value = 0.318009981276496
$ws.Range("T5").Value = 0.318009981276496

# Row 6
$ws.Range("I6").Value = 0.9003255417707673
$ws.Range("J6").Value = 0.9003255417707672
$ws.Range("O6").Value = 0.5533024543641269
$ws.Range("P6").Value = 0.5533024543641269
$ws.Range("S6").Value = 0.4981523319884779
$ws.Range("T6").Value = 0.4981523319884778

# Row 7
$ws.Range("I7").Value = 0.9003255417707673
$ws.Range("J7").Value = 0.9003255417707672
$ws.Range("O7").Value = 0.09348088508103472
$ws.Range("P7").Value = 0.09348088508103473
$ws.Range("S7").Value = 0.08416322850579343
$ws.Range("T7").Value = 0.08416322850579343

# Row 8
$ws.Range("I8").Value = 0.08566479914627706
$ws.Range("J8").Value = 0.08566479914627706
$ws.Range("M8").Value = 97.57717366666668
$ws.Range("N8").Value = 292.731521
$ws.Range("O8").Value = 0.3532166605548384
$ws.Range("P8").Value = 0.3532166605548384
$ws.Range("Q8").Value = 314.8990382621443
$ws.Range("R8").Value = 2834.091344359299
$ws.Range("S8").Value = 0.03025823428154895
$ws.Range("T8").Value = 0.03025823428154896

# Row 9
$ws.Range("I9").Value = 0.08566479914627706
$ws.Range("J9").Value = 0.08566479914627706
$ws.Range("O9").Value = 0.5533024543641269
$ws.Range("P9").Value = 0.5533024543641269
$ws.Range("S9").Value = 0.04739854362024506
$ws.Range("T9").Value = 0.04739854362024506

# Row 10
$ws.Range("I10").Value = 0.08566479914627706
$ws.Range("J10").Value = 0.08566479914627706
$ws.Range("O10").Value = 0.09348088508103472
$ws.Range("P10").Value = 0.09348088508103473
$ws.Range("S10").Value = 0.008008021244483047
$ws.Range("T10").Value = 0.008008021244483049
